$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 103, shifting existing rows 103:127 down to 104:128.
$ws.Rows.Item(103).Insert()

# Populate the newly inserted row 103 with the new record's data.
$ws.Range("A103").Value = 5
$ws.Range("B103").Value = "Macroferia Regional de Talca"
$ws.Range("C103").Value = "Maule"
$ws.Range("D103").Value = 44522
$ws.Range("E103").Value = 7
$ws.Range("F103").Value = 100112024
$ws.Range("G103").Value = "Choclo"
$ws.Range("H103").Value = "Dulce o Americano"
$ws.Range("I103").Value = "Primera"
$ws.Range("J103").Value = 150
$ws.Range("K103").Value = 15000
$ws.Range("L103").Value = 15000
$ws.Range("M103").Value = 15000
$ws.Range("N103").Value = "$/malla 60 unidades"
$ws.Range("O103").Value = "Región de Arica y Parinacota"
$ws.Range("P103").Value = 250
$ws.Range("Q103").Value = 60
$ws.Range("R103").Value = "Hortaliza"
